$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price (column D) and 1h volume-change (column E) figures for each
# coin row, as refreshed by the scheduled GitHub Actions scrape.
$updates = @(
    @{ Row = 2; D = "27.391.99"; E = "  +3.48%  " },
    @{ Row = 3; D = "1.841.51"; E = "  +3.77%  " },
    @{ Row = 4; D = $null; E = "  +2.74%  " },
    @{ Row = 5; D = "318.68"; E = "  +4.59%  " },
    @{ Row = 6; D = "1.024"; E = "  +2.55%  " },
    @{ Row = 7; D = "0.4352"; E = "  +1.67%  " },
    @{ Row = 8; D = "0.3718"; E = "  +2.15%  " },
    @{ Row = 9; D = "0.07341"; E = "  +2.60%  " },
    @{ Row = 10; D = "0.8757"; E = "  +3.33%  " },
    @{ Row = 11; D = "21.36"; E = "  +4.43%  " },
    @{ Row = 12; D = "1.989.34"; E = "  +10.77%  " },
    @{ Row = 13; D = "5.479"; E = "  +4.13%  " },
    @{ Row = 14; D = "6.683"; E = "  +3.52%  " },
    @{ Row = 15; D = "0.07163"; E = "  +4.34%  " },
    @{ Row = 16; D = "82.20"; E = "  +4.18%  " },
    @{ Row = 17; D = $null; E = "  +2.49%  " },
    @{ Row = 18; D = "0.000008985"; E = "  +3.38%  " },
    @{ Row = 19; D = "1.022"; E = "  +2.44%  " },
    @{ Row = 20; D = "15.40"; E = "  +2.51%  " },
    @{ Row = 21; D = "27.424.44"; E = "  +3.55%  " },
    @{ Row = 22; D = "5.259"; E = "  +2.87%  " },
    @{ Row = 23; D = "11.13"; E = "  +0.20%  " },
    @{ Row = 24; D = "2.191.03"; E = "  +9.83%  " },
    @{ Row = 25; D = "156.86"; E = "  +3.03%  " },
    @{ Row = 26; D = "1.904"; E = "  +1.76%  " },
    @{ Row = 27; D = "18.55"; E = "  +2.96%  " },
    @{ Row = 28; D = "5.276"; E = "  +3.80%  " },
    @{ Row = 29; D = "1.923"; E = "  +6.68%  " },
    @{ Row = 30; D = "115.50"; E = "  +1.44%  " },
    @{ Row = 31; D = "0.09016"; E = "  +0.95%  " },
    @{ Row = 32; D = "1.200"; E = "  +6.27%  " },
    @{ Row = 33; D = "0.7598"; E = "  +4.29%  " },
    @{ Row = 34; D = "4.474"; E = "  +3.30%  " },
    @{ Row = 35; D = "2.858"; E = "  +4.29%  " },
    @{ Row = 36; D = "1.026"; E = "  +2.77%  " },
    @{ Row = 37; D = "1.149"; E = "  +3.62%  " },
    @{ Row = 38; D = $null; E = "  +3.64%  " },
    @{ Row = 39; D = "0.05251"; E = "  +1.82%  " },
    @{ Row = 40; D = "0.5155"; E = "  +4.38%  " },
    @{ Row = 41; D = "2.801"; E = "  +6.72%  " },
    @{ Row = 42; D = "0.1661"; E = "  +2.95%  " },
    @{ Row = 43; D = "6.532"; E = "  +3.40%  " },
    @{ Row = 44; D = "8.470"; E = "  +5.65%  " },
    @{ Row = 45; D = "108.31"; E = "  +3.17%  " },
    @{ Row = 46; D = "10.54"; E = "  +3.59%  " },
    @{ Row = 47; D = $null; E = "  +2.89%  " },
    @{ Row = 48; D = "0.4632"; E = "  +3.45%  " },
    @{ Row = 49; D = "1.670"; E = "  +2.21%  " },
    @{ Row = 50; D = $null; E = "  +9.04%  " },
    @{ Row = 51; D = "0.06291"; E = "  +1.39%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D" + $u.Row)
        # Force text storage so values such as "1.024" or "27.391.99" are not
        # reinterpreted as numbers/dates by Excel, matching the source data
        # which stores these as plain text strings.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    $ws.Range("E" + $u.Row).Value = $u.E
}
